$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 90, shifting existing row 90 (and below) down to row 91.
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the new record.
$ws.Range("A90").Value = 10
$ws.Range("B90").Value = "Vega Modelo de Temuco"
$ws.Range("C90").Value = "La Araucanía"
$ws.Range("D90").Value = 45229
$ws.Range("E90").Value = 9
$ws.Range("F90").Value = 100112022
$ws.Range("G90").Value = "Arveja Verde"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 185
$ws.Range("K90").Value = 24000
$ws.Range("L90").Value = 24000
$ws.Range("M90").Value = 24000
$ws.Range("N90").Value = "$/saco 25 kilos"
$ws.Range("O90").Value = "Región del Maule"
$ws.Range("P90").Value = 960
$ws.Range("Q90").Value = 25
$ws.Range("R90").Value = "Hortaliza"
